$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, shifting existing rows 25-126 down to 26-127
$ws.Rows.Item(25).EntireRow.Insert()

# Populate the newly inserted row 25 with data
$ws.Range("A25").Value = 1
$ws.Range("B25").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C25").Value = "Arica y Parinacota"
$ws.Range("D25").Value = 44811
$ws.Range("D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E25").Value = 15
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100106
$ws.Range("H25").Value = "Oleaginosos"
$ws.Range("I25").Value = 100106002
$ws.Range("J25").Value = "Palta"
$ws.Range("K25").Value = "Edranol"
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 59000
$ws.Range("O25").Value = 60000
$ws.Range("P25").Value = 59500
$ws.Range("Q25").Value = "`$/caja 25 kilos"
$ws.Range("R25").Value = "Región de Coquimbo"
$ws.Range("S25").Value = 2380
$ws.Range("T25").Value = 25
